# Update market-price / profit figures pulled by the scheduled pricing runner.
# Mirrors xl/worksheets/sheet*.xml cell-value changes (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H:N) captured in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 365.16666
$ws.Range("I6").Value2 = 295.5
$ws.Range("K6").Value2 = 886.5
$ws.Range("M6").Value2 = -774.5

$ws.Range("H41").Value2 = 495
$ws.Range("J41").Value2 = 616.5
$ws.Range("L41").Value2 = 616.5
$ws.Range("N41").Value2 = -1496.5

$ws.Range("H116").Value2 = 4538.077
$ws.Range("I116").Value2 = 3695.375
$ws.Range("J116").Value2 = 5886.4
$ws.Range("K116").Value2 = 3695.375
$ws.Range("L116").Value2 = 5886.4
$ws.Range("M116").Value2 = -253.375
$ws.Range("N116").Value2 = -12770.4

$ws.Range("H132").Value2 = 5988.5
$ws.Range("I132").Value2 = 5988
$ws.Range("K132").Value2 = 17964
$ws.Range("M132").Value2 = -15434

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 888.75
$ws.Range("I4").Value2 = 888.75
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 888.75
$ws.Range("L4").Value2 = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value2 = -772.75

$ws.Range("H32").Value2 = 13956.667
$ws.Range("I32").Value2 = 12498.182
$ws.Range("K32").Value2 = 12498.182
$ws.Range("M32").Value2 = -12211.182

$ws.Range("H94").Value2 = 73832.25
$ws.Range("J94").Value2 = 73832.25
$ws.Range("L94").Value2 = 73832.25
$ws.Range("N94").Value2 = -75634.25

$ws.Range("H110").Value2 = 4282.6665
$ws.Range("I110").Value2 = 696.5
$ws.Range("K110").Value2 = 696.5
$ws.Range("M110").Value2 = 1348.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 130
$ws.Range("I22").Value2 = 130
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 130
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value2 = 43
$ws.Range("N22").Value2 = 0

$ws.Range("H99").Value2 = 3099.8572
$ws.Range("I99").Value2 = 3099.8572
$ws.Range("K99").Value2 = 3099.8572
$ws.Range("M99").Value2 = -1601.8572

$ws.Range("H134").Value2 = 9166.200000000001
$ws.Range("I134").Value2 = 10449.5
$ws.Range("J134").Value2 = 6599.6
$ws.Range("K134").Value2 = 31348.5
$ws.Range("L134").Value2 = 19798.8
$ws.Range("M134").Value2 = -28813.5
$ws.Range("N134").Value2 = -24868.8

$ws.Range("H141").Value2 = 30000
$ws.Range("I141").Value2 = 30000
$ws.Range("J141").Value2 = 0
$ws.Range("K141").Value2 = 30000
$ws.Range("L141").ClearContents()
$ws.Range("M141").Value2 = -24820
$ws.Range("N141").Value2 = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 66.25
$ws.Range("I7").Value2 = 75.07692
$ws.Range("J7").Value2 = 49.857143
$ws.Range("K7").Value2 = 75.07692
$ws.Range("L7").Value2 = 49.857143
$ws.Range("M7").Value2 = 37.92308
$ws.Range("N7").Value2 = -275.857143

$ws.Range("H41").Value2 = 15699.714
$ws.Range("J41").Value2 = 21250
$ws.Range("L41").Value2 = 21250
$ws.Range("N41").Value2 = -22106

$ws.Range("H58").Value2 = 2122.2856
$ws.Range("I58").Value2 = 2098.889
$ws.Range("J58").Value2 = 2262.6667
$ws.Range("K58").Value2 = 2098.889
$ws.Range("L58").Value2 = 2262.6667
$ws.Range("M58").Value2 = -1895.889
$ws.Range("N58").Value2 = -2668.6667

$ws.Range("H99").Value2 = 5234.7
$ws.Range("I99").Value2 = 6068.5
$ws.Range("J99").Value2 = 1899.5
$ws.Range("K99").Value2 = 6068.5
$ws.Range("L99").Value2 = 1899.5
$ws.Range("M99").Value2 = -4570.5
$ws.Range("N99").Value2 = -4895.5

$ws.Range("H126").Value2 = 5234.7
$ws.Range("I126").Value2 = 6068.5
$ws.Range("J126").Value2 = 1899.5
$ws.Range("K126").Value2 = 18205.5
$ws.Range("L126").Value2 = 5698.5
$ws.Range("M126").Value2 = -15735.5
$ws.Range("N126").Value2 = -10638.5

$ws.Range("H132").Value2 = 1499.5
$ws.Range("I132").Value2 = 1499.5
$ws.Range("K132").Value2 = 4498.5
$ws.Range("M132").Value2 = -1968.5

$ws.Range("H136").Value2 = 2122.2856
$ws.Range("I136").Value2 = 2098.889
$ws.Range("J136").Value2 = 2262.6667
$ws.Range("K136").Value2 = 6296.667
$ws.Range("L136").Value2 = 6788.000100000001
$ws.Range("M136").Value2 = -3746.667
$ws.Range("N136").Value2 = -11888.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 1182.5
$ws.Range("I131").Value2 = 1182.5
$ws.Range("K131").Value2 = 3547.5
$ws.Range("M131").Value2 = 1492.5

$ws.Range("H137").Value2 = 3103.5557
$ws.Range("I137").Value2 = 2499.8
$ws.Range("J137").Value2 = 3858.25
$ws.Range("K137").Value2 = 7499.400000000001
$ws.Range("L137").Value2 = 11574.75
$ws.Range("M137").Value2 = -2399.400000000001
$ws.Range("N137").Value2 = -21774.75

$ws.Range("H138").Value2 = 6128.875
$ws.Range("I138").Value2 = 3676.6667
$ws.Range("K138").Value2 = 11030.0001
$ws.Range("M138").Value2 = -5890.000100000001

$ws.Range("H139").Value2 = 4578.909
$ws.Range("I139").Value2 = 3374.2222
$ws.Range("K139").Value2 = 10122.6666
$ws.Range("M139").Value2 = -4982.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 435.2
$ws.Range("I2").Value2 = 435.2
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 435.2
$ws.Range("L2").Value2 = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value2 = -322.2

$ws.Range("I80").Value2 = 2085
$ws.Range("J80").Value2 = 1953.76
$ws.Range("K80").Value2 = 2085
$ws.Range("L80").Value2 = 1953.76
$ws.Range("M80").Value2 = -1087
$ws.Range("N80").Value2 = -3949.76

$ws.Range("I83").Value2 = 2085
$ws.Range("J83").Value2 = 1953.76
$ws.Range("K83").Value2 = 10425
$ws.Range("L83").Value2 = 9768.799999999999
$ws.Range("M83").Value2 = -5433
$ws.Range("N83").Value2 = -19752.8

$ws.Range("H113").Value2 = 1933.3334

$ws.Range("H126").Value2 = 2832.1667
$ws.Range("I126").Value2 = 2798.6
$ws.Range("J126").Value2 = 3000
$ws.Range("K126").Value2 = 8395.799999999999
$ws.Range("L126").Value2 = 9000
$ws.Range("M126").Value2 = -5925.799999999999
$ws.Range("N126").Value2 = -13940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value2 = 0
$ws.Range("J18").Value2 = 0
$ws.Range("L18").ClearContents()
$ws.Range("N18").Value2 = 0

$ws.Range("H61").Value2 = 8024.778
$ws.Range("I61").Value2 = 7869.6665
$ws.Range("K61").Value2 = 7869.6665
$ws.Range("M61").Value2 = -7667.6665

$ws.Range("H105").Value2 = 0
$ws.Range("J105").Value2 = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value2 = 0

$ws.Range("H113").Value2 = 8024.778
$ws.Range("I113").Value2 = 7869.6665
$ws.Range("K113").Value2 = 7869.6665
$ws.Range("M113").Value2 = -5699.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value2 = 2444.4443
$ws.Range("J2").Value2 = 2000
$ws.Range("L2").Value2 = 2000
$ws.Range("N2").Value2 = -2224

$ws.Range("H81").Value2 = 300
$ws.Range("J81").Value2 = 0
$ws.Range("L81").Value2 = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value2 = 300
$ws.Range("J84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("N84").ClearContents()

$ws.Range("H96").Value2 = 712.5
$ws.Range("I96").Value2 = 712.5
$ws.Range("J96").Value2 = 0
$ws.Range("K96").Value2 = 712.5
$ws.Range("L96").Value2 = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value2 = 660.5

$ws.Range("H136").Value2 = 1757.3462
$ws.Range("I136").Value2 = 1508.3043
$ws.Range("J136").Value2 = 3666.6667
$ws.Range("K136").Value2 = 4524.9129
$ws.Range("L136").Value2 = 11000.0001
$ws.Range("M136").Value2 = -1974.9129
$ws.Range("N136").Value2 = -16100.0001
